# 2021-05 Victorian Outbreak Paths.xlsx -- "Add files via upload"
#
# Sheet1 ("Sheet1"): one new data row (110) appended to Table1, continuing
#   the outbreak-path log with a new Arcare Maidstone link.
# Sheet2 ("Date Colours"): the purple date-gradient (column B) is
#   recomputed to smoothly span the now-larger date range (34 -> 38 rows)
#   and four new trailing dates (35-38) are appended to Table "Date_Colours".
# The "Date Colours" tab becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: append the new outbreak-path row
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$lo1 = $ws1.ListObjects.Item("Table1")

$newRow1 = $lo1.ListRows.Add()
$r1 = $newRow1.Range.Row

$ws1.Cells.Item($r1, 1).Value = 44377
$ws1.Cells.Item($r1, 1).NumberFormat = "d-mmm"
$ws1.Cells.Item($r1, 2).Value = "A # l"
$ws1.Cells.Item($r1, 3).Value = "A # m"
$ws1.Cells.Item($r1, 4).Value = "Arcare Maidstone"
$ws1.Cells.Item($r1, 6).Value = "Epping private hospital"
$ws1.Cells.Item($r1, 7).Value = "Kappa (B.1.617.1)"

$ws1.Range("D" + $r1).Select()

# ---------------------------------------------------------------------
# Sheet2 ("Date Colours"): refresh the purple gradient for the existing
# rows, then append four new rows for the latest dates.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Date Colours")
$lo2 = $ws2.ListObjects.Item("Date_Colours")

$gradient = @(
  "#fbfbff","#f7f6ff","#f3f2fe","#efedfe","#ebe9fe","#e7e5fe","#e3e0fd","#dfdcfd",
  "#dbd7fd","#d7d3fd","#d3cffc","#cfcbfc","#cbc6fc","#c6c2fb","#c2befb","#bebafb",
  "#b9b5fa","#b5b1fa","#b1adfa","#aca9f9","#a7a5f9","#a3a1f8","#9e9df8","#9999f7",
  "#9495f7","#8f91f7","#8a8df6","#8589f6","#7f85f5","#7a81f5","#747df4","#6e79f4",
  "#6775f3","#6171f3","#596df2","#526af2","#4966f1"
)

# Existing rows 2..34 -> gradient[0..32]; update the "Colour Code" column only.
for ($i = 0; $i -lt 33; $i++) {
  $row = 2 + $i
  $ws2.Cells.Item($row, 2).Value = $gradient[$i]
}

# Four brand-new dates (21-Jun through 24-Jun 2021) appended to the table.
$newDates = @(44374, 44375, 44376, 44377)
for ($j = 0; $j -lt 4; $j++) {
  $newRow2 = $lo2.ListRows.Add()
  $r2 = $newRow2.Range.Row
  $ws2.Cells.Item($r2, 1).Value = $newDates[$j]
  $ws2.Cells.Item($r2, 1).NumberFormat = "d-mmm"
  $ws2.Cells.Item($r2, 2).Value = $gradient[33 + $j]
  $ws2.Cells.Item($r2, 4).Value = "Diamond"
}

# The final new row also carries the day's highlight colour.
$ws2.Cells.Item($r2, 3).Value = "#EF7C34"

$ws2.Range("D34:D" + $r2).Select()

# "Date Colours" becomes the active sheet/tab.
$ws2.Activate()
